# "Finished week 2 reading" - fill in actual time spent reading for week 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week2")

$ws.Range("C4").Value = 30/1440
$ws.Range("C5").Value = 15/1440
$ws.Range("C6").Value = 25/1440
$ws.Range("C7").Value = 6/1440

# Move the active cell selection on this sheet to C5 (as in the source file)
$ws.Activate()
$ws.Range("C5").Select()
